$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.685.51'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.69%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.166.17'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.80%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '615.82'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.86%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.89'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.67%  '
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.165.31'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.77%  '
$ws.Range("E9").Value = '  -1.03%  '
$ws.Range("E10").Value = '  -1.13%  '
$ws.Range("E11").Value = '  -2.32%  '
$ws.Range("E12").Value = '  -1.30%  '
$ws.Range("E13").Value = '  -0.32%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.87'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.30%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.687.53'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.02%  '
$ws.Range("E16").Value = '  +2.90%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.721.51'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.58%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.165.52'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.96%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.93'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.83%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '479.17'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.07%  '
$ws.Range("E21").Value = '  -0.24%  '
$ws.Range("E22").Value = '  +1.34%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.96'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.54%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.81'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.89%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.59'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.64%  '
$ws.Range("E26").Value = '  -0.07%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.83'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.14%  '
$ws.Range("E28").Value = '  -1.22%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.94'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.56%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.118'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -6.36%  '
$ws.Range("E31").Value = '  -7.94%  '
$ws.Range("E32").Value = '  +0.10%  '
$ws.Range("E33").Value = '  -0.33%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '26.57'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.28%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.13'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.06%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0₃0778'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.34%  '
$ws.Range("E37").Value = '  -1.94%  '
$ws.Range("B38").Value = 'dogwifhat'
$ws.Range("C38").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.19'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.17%  '
$ws.Range("B39").Value = 'OKB'
$ws.Range("C39").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '53.05'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.03%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '461.54'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.16%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0399'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.10%  '
$ws.Range("E42").Value = '  -4.22%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.42'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.60%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.851.47'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.66%  '
$ws.Range("E45").Value = '  -3.85%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.269'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.49%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.43'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.01%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '26.67'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.59%  '
$ws.Range("E49").Value = '  +0.08%  '
$ws.Range("E50").Value = '  -1.73%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '120.48'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.21%  '
